$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.195.23'
$ws.Range("E2").Value = '  -2.93%  '

$ws.Range("D3").Value = '1.609.42'
$ws.Range("E3").Value = '  -2.38%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9998'
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.82'
$ws.Range("E6").Value = '  -2.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3781'
$ws.Range("E7").Value = '  -2.94%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3667'
$ws.Range("E8").Value = '  -4.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.45'
$ws.Range("E9").Value = '  -4.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9985'
$ws.Range("E10").Value = '  -0.27%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.267'
$ws.Range("E11").Value = '  -6.25%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08085'
$ws.Range("E12").Value = '  -4.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.98'
$ws.Range("E13").Value = '  -3.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.604'
$ws.Range("E14").Value = '  -6.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.424'
$ws.Range("E15").Value = '  -6.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001255'
$ws.Range("E16").Value = '  -4.66%  '

$ws.Range("D17").Value = '1.610.00'
$ws.Range("E17").Value = '  -2.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.42'
$ws.Range("E18").Value = '  -3.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06790'
$ws.Range("E19").Value = '  -2.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.38'
$ws.Range("E20").Value = '  -6.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.577'
$ws.Range("E21").Value = '  -5.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("E23").Value = '  -4.77%  '

$ws.Range("D24").Value = '23.211.98'
$ws.Range("E24").Value = '  -2.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.352'
$ws.Range("E25").Value = '  -4.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.840'
$ws.Range("E26").Value = '  -3.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.09'
$ws.Range("E27").Value = '  -4.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.46'
$ws.Range("E28").Value = '  -0.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.290'
$ws.Range("E29").Value = '  -2.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.56'
$ws.Range("E30").Value = '  -4.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.405'
$ws.Range("E31").Value = '  -4.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.863'
$ws.Range("E32").Value = '  -12.76%  '

$ws.Range("D33").Value = '1.789.04'
$ws.Range("E33").Value = '  -2.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9697'
$ws.Range("E34").Value = '  -7.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07699'
$ws.Range("E35").Value = '  -4.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02768'
$ws.Range("E36").Value = '  -6.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2560'
$ws.Range("E37").Value = '  -4.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.245'
$ws.Range("E38").Value = '  -6.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.19'
$ws.Range("E39").Value = '  -7.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08904'
$ws.Range("E40").Value = '  -2.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.394'
$ws.Range("E41").Value = '  -1.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7196'
$ws.Range("E42").Value = '  -5.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.82'
$ws.Range("E43").Value = '  -4.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.71'
$ws.Range("E44").Value = '  -2.89%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6658'
$ws.Range("E45").Value = '  -4.70%  '

$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.307'
$ws.Range("E46").Value = '  -6.40%  '

$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9992'
$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("E48").Value = '  -2.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08005'
$ws.Range("E49").Value = '  -3.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.85'
$ws.Range("E50").Value = '  -2.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.179'
$ws.Range("E51").Value = '  -3.01%  '

